# Rename the two existing sheets, add a third ("156") as the new active
# sheet, and build out its "Phiếu nhập kho / Phiếu xuất kho" (stock
# in/out slip) form layout, matching the commit "queue and queue manager
# should be done".

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "idea"
$wb.Worksheets.Item(2).Name = "shortVersion"

# Add the new sheet after the last existing sheet so it lands at the end
# and becomes the active tab, mirroring the target workbook.xml
# (activeTab="2", tabSelected on the new sheet).
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "156"

# --- "Phiếu nhập kho" (stock-in slip) block ------------------------------
$ws.Range("C5").Value = "Phiếu nhập kho"
$ws.Range("C5").Font.Color = 255

$ws.Range("D6").Value = "Tên người lập"
$ws.Range("F6").Value = "userName"
$ws.Range("I6").Value = "Ngày lập"
$ws.Range("J6").Value = "(Tạm bỏ qua)"

$ws.Range("D7").Value = "Tên sản phẩm"
$ws.Range("F7").Value = "Nhập tay"
$ws.Range("I7").Value = "ID phiếu"
$ws.Range("J7").Value = "Auto"

$ws.Range("D8").Value = "Số lượng"
$ws.Range("F8").Value = "Nhập tay"
$ws.Range("I8").Value = "Lý do nhập"
$ws.Range("J8").Value = "Nhập tay"

$ws.Range("D9").Value = "giá nhập"
$ws.Range("F9").Value = "Nhập tay"
$ws.Range("N9").Value = "Tăng 156"
$ws.Range("O9").Value = "Giảm 111"
$ws.Range("P9").Value = "và"
$ws.Range("Q9").Value = "tăng 632"

$ws.Range("B10").Value = ">>>"
$ws.Range("C10").Value = "Set giá bán sau khi tạo phiếu nhập kho"
$ws.Range("G10").Value = "Nhập tay"
$ws.Range("L10").Value = ">>> Lập bút toán"
$ws.Range("N10").Value = "số lượng x giá nhập"

# Row 11 is a tall spacer row beneath the "nhập kho" block.
$ws.Rows.Item(11).RowHeight = 45

# --- "Phiếu xuất kho" (stock-out slip) block ------------------------------
$ws.Range("C12").Value = "Phiếu xuất kho"
$ws.Range("C12").Font.Color = 255

$ws.Range("D13").Value = "Tên người lập"
$ws.Range("F13").Value = "userName"
$ws.Range("I13").Value = "Ngày lập"
$ws.Range("J13").Value = "(Tạm bỏ qua)"

$ws.Range("D14").Value = "Tên sản phẩm"
$ws.Range("F14").Value = "Nhập tay"
$ws.Range("I14").Value = "ID phiếu"
$ws.Range("J14").Value = "Auto"

$ws.Range("D15").Value = "Số lượng"
$ws.Range("F15").Value = "Nhập tay"
$ws.Range("I15").Value = "Lý do xuất"
$ws.Range("J15").Value = "Nhập tay"

$ws.Range("C16").Value = "switchcase>>>>"
$ws.Range("D16").Value = "giá bán"
$ws.Range("F16").Value = "Lấy từ giá set ở trên"
$ws.Range("L16").Value = ">>> Lập bút toán"
$ws.Range("N16").Value = "giảm 156"
$ws.Range("O16").Value = "tăng 111"
$ws.Range("P16").Value = "và"
$ws.Range("Q16").Value = "tăng 515 "

$ws.Range("D17").Value = "hoặc giá bán = 0"
$ws.Range("F17").Value = "hàng hỏng/ lý do khác"
$ws.Range("N17").Value = "tính theo phương pháp fifo"

# --- column widths ---------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 19.833333333333336
$ws.Columns.Item(6).ColumnWidth = 19.333333333333336
$ws.Columns.Item(14).ColumnWidth = 25.166666666666668
$ws.Columns.Item(15).ColumnWidth = 13.166666666666668

# --- page setup / view state -------------------------------------------
$ws.PageSetup.Orientation = 1

$ws.Range("H11").Select() | Out-Null
